$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.77"
$ws.Range("E2").Value = "'5.28%"
$ws.Range("D3").Value = "'44.48"
$ws.Range("E3").Value = "'7.05%"
$ws.Range("D4").Value = "'5.103"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.08003"
$ws.Range("E5").Value = "'5.97%"
$ws.Range("D6").Value = "'4.491"
$ws.Range("E6").Value = "'2.65%"
$ws.Range("D7").Value = "'1.645"
$ws.Range("E7").Value = "'2.97%"
$ws.Range("D8").Value = "'1.085"
$ws.Range("E8").Value = "'16.95%"
$ws.Range("D9").Value = "'0.1293"
$ws.Range("E9").Value = "'7.07%"
$ws.Range("D10").Value = "'0.1890"
$ws.Range("E10").Value = "'3.05%"
$ws.Range("D11").Value = "'0.09291"
$ws.Range("E11").Value = "'4.72%"
$ws.Range("D12").Value = "'0.04200"
$ws.Range("E12").Value = "'7.29%"
$ws.Range("D13").Value = "'0.1040"
$ws.Range("E13").Value = "'-1.13%"
$ws.Range("D14").Value = "'0.001305"
$ws.Range("E14").Value = "'2.17%"
$ws.Range("D15").Value = "'0.005856"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("D17").Value = "'3.373"
$ws.Range("E17").Value = "'1.20%"
$ws.Range("E18").Value = "'-0.79%"
$ws.Range("E19").Value = "'1.23%"
$ws.Range("D20").Value = "'8.021"
$ws.Range("E20").Value = "'1.07%"
$ws.Range("D21").Value = "'0.1348"
$ws.Range("E21").Value = "'-4.93%"
$ws.Range("D23").Value = "'0.04195"
$ws.Range("E23").Value = "'3.32%"
$ws.Range("D24").Value = "'0.001272"
$ws.Range("E24").Value = "'0.61%"
$ws.Range("D25").Value = "'0.004591"
$ws.Range("E25").Value = "'14.94%"
$ws.Range("D26").Value = "'0.0001339"
$ws.Range("E26").Value = "'8.85%"
$ws.Range("D38").Value = "'0.02650"
$ws.Range("E38").Value = "'9.74%"
$ws.Range("D39").Value = "'0.05423"
$ws.Range("E39").Value = "'4.16%"
$ws.Range("D40").Value = "'0.005618"
$ws.Range("E40").Value = "'-12.06%"
$ws.Range("D41").Value = "'0.007730"
$ws.Range("E41").Value = "'-0.78%"
$ws.Range("D42").Value = "'0.1416"
$ws.Range("E42").Value = "'6.64%"
$ws.Range("D43").Value = "'0.007326"
$ws.Range("E43").Value = "'-3.30%"
$ws.Range("D44").Value = "'0.008603"
$ws.Range("E44").Value = "'9.90%"
$ws.Range("D45").Value = "'0.3115"
$ws.Range("E45").Value = "'-3.40%"
$ws.Range("D46").Value = "'0.00006740"
$ws.Range("E46").Value = "'-0.54%"
$ws.Range("D47").Value = "'0.00000000743"
$ws.Range("E47").Value = "'-0.85%"
$ws.Range("D48").Value = "'0.05457"
$ws.Range("E48").Value = "'21.21%"
$ws.Range("D49").Value = "'0.003963"
$ws.Range("E49").Value = "'-5.57%"
$ws.Range("D50").Value = "'0.00002081"
$ws.Range("E50").Value = "'-0.85%"
$ws.Range("D51").Value = "'0.0001982"
$ws.Range("E51").Value = "'-0.85%"
